$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.209.93"
$ws.Range("E2").Value = "  -0.59%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.783.48"
$ws.Range("E3").Value = "  -1.78%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.89"
$ws.Range("E5").Value = "  -2.38%  "

$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3777"
$ws.Range("E7").Value = "  -2.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3432"
$ws.Range("E8").Value = "  -3.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.32"
$ws.Range("E9").Value = "  -4.15%  "

$ws.Range("E10").Value = "  -4.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07480"
$ws.Range("E11").Value = "  -4.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.70"
$ws.Range("E13").Value = "  -4.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.461"
$ws.Range("E14").Value = "  -3.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.786.20"
$ws.Range("E15").Value = "  -1.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.091"
$ws.Range("E16").Value = "  -2.58%  "

$ws.Range("E17").Value = "  -3.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06665"
$ws.Range("E18").Value = "  -1.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.83"
$ws.Range("E19").Value = "  -3.84%  "

$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.610"
$ws.Range("E21").Value = "  -0.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.31"
$ws.Range("E22").Value = "  -4.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.216.35"
$ws.Range("E23").Value = "  -0.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.38"
$ws.Range("E24").Value = "  -6.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.411"
$ws.Range("E25").Value = "  -2.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.506"
$ws.Range("E26").Value = "  -0.94%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.540"
$ws.Range("E27").Value = "  -7.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.32"
$ws.Range("E28").Value = "  -3.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.61"
$ws.Range("E29").Value = "  -0.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.989.63"
$ws.Range("E30").Value = "  -1.29%  "

$ws.Range("E31").Value = "  -2.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.017"
$ws.Range("E32").Value = "  -2.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.080"
$ws.Range("E33").Value = "  -6.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08687"
$ws.Range("E34").Value = "  -2.01%  "

$ws.Range("E35").Value = "  -5.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.660"
$ws.Range("E36").Value = "  -3.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6946"
$ws.Range("E37").Value = "  -1.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.458"
$ws.Range("E38").Value = "  -4.27%  "

$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2199"
$ws.Range("E39").Value = "  -3.64%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06331"
$ws.Range("E40").Value = "  -4.04%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.792"
$ws.Range("E41").Value = "  -3.11%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02338"
$ws.Range("E42").Value = "  -4.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.244"
$ws.Range("E43").Value = "  -1.65%  "

$ws.Range("E44").Value = "  -4.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6499"
$ws.Range("E45").Value = "  -2.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.851"
$ws.Range("E47").Value = "  -2.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.145"
$ws.Range("E48").Value = "  -2.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.12"
$ws.Range("E49").Value = "  -3.48%  "

$ws.Range("E50").Value = "  -3.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.21"
$ws.Range("E51").Value = "  -2.60%  "
